# Apply the edits described by the diff:
#  - Add a "Note:" row (row 11) on the Setup sheet: a label cell and a
#    highlighted (yellow-filled) note cell.
#  - Switch the active/selected sheet from "Copy To Nodes" back to "Setup".
#  - Update the selection on each sheet to match the new saved view state.

$wb = $excel.ActiveWorkbook
$wsSetup = $wb.Worksheets.Item("Setup")
$wsCopy  = $wb.Worksheets.Item("Copy To Nodes")

# --- Add the new "Note" row on the Setup sheet ---
$wsSetup.Range("F11").Value = "Note:"
$wsSetup.Range("F11").Style = "Accent1"

$wsSetup.Range("G11").Value = "Not a Microsoft supported deployment topology"
$wsSetup.Range("G11").Interior.Color = 65535

# --- Set/refresh the selection on "Copy To Nodes" (stays A15) ---
$wsCopy.Activate()
$wsCopy.Range("A15").Select()

# --- Make "Setup" the active sheet again (was "Copy To Nodes") and set its selection ---
$wsSetup.Activate()
$wsSetup.Range("G16").Select()
